$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2.8
$ws.Range("J2").Value = 2.67
$ws.Range("K2").Value = 1.8
$ws.Range("N2").Value = 4.3
$ws.Range("O2").Value = 1.6
$ws.Range("P2").Value = 2.07
$ws.Range("T2").Value = 2.05
$ws.Range("W2").Value = 4.6
$ws.Range("Z2").Value = 17
$ws.Range("AB2").Value = 55
$ws.Range("AC2").Value = 4.6
$ws.Range("AD2").Value = 5.9
$ws.Range("AM2").Value = 80
$ws.Range("AN2").Value = 3.5
$ws.Range("AO2").Value = 10.75
$ws.Range("AP2").Value = 27
$ws.Range("AR2").Value = 120
$ws.Range("AT2").Value = 2.02
